$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update header row (row 1) for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 (CON) for columns B:E
$ws.Range("B2").Value = 169.76763831018951
$ws.Range("C2").Value = 136.52373354138246
$ws.Range("D2").Value = 171.04508440339703
$ws.Range("E2").Value = 134.77546343603575

# Update row 3 (STR) for columns B:E
$ws.Range("B3").Value = 153.98837884296645
$ws.Range("C3").Value = 128.60867508149147
$ws.Range("D3").Value = 158.2887464467814
$ws.Range("E3").Value = 129.47784913000453

# Update the selection on the sheet to match new sqref B1:E3
$ws.Range("B1:E3").Select()
